# Auto-emailsender update: record when the first and second reminder
# emails were sent for each asset owner.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns (bold, matching the rest of row 1) -------------
$ws.Range("H1").Value = "First Email"
$ws.Range("H1").Font.Bold = $true

$ws.Range("I1").Value = "Second Email"
$ws.Range("I1").Font.Bold = $true

# --- Timestamp columns (date-only, formatted as yyyy-mm-dd) -----------
# Column I ("Second Email") is populated first, then column H
# ("First Email") - matches the order the sender script touched cells in.
$ws.Range("I2").Value = 43562
$ws.Range("I2").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("I3").Value = 43562
$ws.Range("I3").NumberFormat = "yyyy-mm-dd"

$ws.Range("I4").Value = 43562
$ws.Range("I4").NumberFormat = "yyyy-mm-dd"

$ws.Range("H2").Value = 43562
$ws.Range("H2").NumberFormat = "yyyy-mm-dd"

$ws.Range("H3").Value = 43562
$ws.Range("H3").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("H4").Value = 43562
$ws.Range("H4").NumberFormat = "yyyy\-mm\-dd"

# --- Column widths for the new columns ---------------------------------
$ws.Columns.Item(8).ColumnWidth = 23.57
$ws.Columns.Item(9).ColumnWidth = 13.71

# --- Hyperlinks: B2 loses its hyperlink (already emailed), B3 gets its
#     own individual hyperlink (was merged B2:B3 before), B4 unchanged. --
$ws.Hyperlinks.Item(2).Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:jake.p.burgess@gmail.com") | Out-Null
$ws.Range("B3").Style = "Hyperlink"

# --- Restore selection state --------------------------------------------
$ws.Range("I13").Select() | Out-Null
